$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp banner in A1.
$ws.Range("A1").Value = "Datos actualizados a 12 de Julio de 2020 a las 10:26"

# Swap Eslovaquia / Zambia: Eslovaquia's case count (1901) overtook
# Zambia's (1895), so the two countries trade places in the ranking.
$ws.Range("A118").Value = "Eslovaquia"
$ws.Range("A119").Value = "Zambia"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 3355781
$ws.Range("C4").Value = 135
$ws.Range("D4").Value = 1490542
$ws.Range("E4").Value = 1727836

# Row 7 - Rusia
$ws.Range("B7").Value = 727162
$ws.Range("C7").Value = 6615
$ws.Range("D7").Value = 501061
$ws.Range("E7").Value = 214766
$ws.Range("G7").Value = 130
$ws.Range("H7").Value = 11335

# Row 39 - Ucrania
$ws.Range("B39").Value = 53521
$ws.Range("C39").Value = 678
$ws.Range("D39").Value = 26118
$ws.Range("E39").Value = 26020
$ws.Range("G39").Value = 11
$ws.Range("H39").Value = 1383

# Row 43 - Singapur
$ws.Range("B43").Value = 45961
$ws.Range("C43").Value = 178
$ws.Range("E43").Value = 3909

# Row 60 - Moldavia
$ws.Range("D60").Value = 12667
$ws.Range("E60").Value = 5900
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 641

# Row 61 - Austria
$ws.Range("B61").Value = 18897
$ws.Range("C61").Value = 114
$ws.Range("D61").Value = 16952
$ws.Range("E61").Value = 1239

# Row 76 - El Salvador
$ws.Range("E76").Value = 3580
$ws.Range("G76").Value = 6
$ws.Range("H76").Value = 260

# Row 118 - now Eslovaquia (was Zambia's row by position, ranking moved up)
$ws.Range("B118").Value = 1901
$ws.Range("C118").Value = 8
$ws.Range("D118").Value = 1493
$ws.Range("E118").Value = 380
$ws.Range("H118").Value = 28

# Row 119 - now Zambia
$ws.Range("B119").Value = 1895
$ws.Range("C119").Value = 0
$ws.Range("D119").Value = 1348
$ws.Range("E119").Value = 505
$ws.Range("H119").Value = 42

# Row 123 - Eslovenia
$ws.Range("B123").Value = 1841
$ws.Range("C123").Value = 14
$ws.Range("E123").Value = 301
